# Update the "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# values on the zh-cn and de-de report sheets (row 5 - the f34aada8... entry)
# to reflect a freshly generated handback report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-01-27 08:46:21"
$wsZhCn.Range("G5").Value = "2016-01-27 08:47:08"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-01-27 08:46:33"
$wsDeDe.Range("G5").Value = "2016-01-27 08:47:28"
